$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 307, shifting existing rows 307:327 down to 308:328
$ws.Rows.Item(307).Insert()

# Fill in the new row 307 (copy of the former row 307 data with updated values)
$ws.Cells.Item(307, 1).Value = 5
$ws.Cells.Item(307, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(307, 3).Value = "Maule"
$ws.Cells.Item(307, 4).Value = 44714
$ws.Cells.Item(307, 5).Value = 7
$ws.Cells.Item(307, 6).Value = 100114014
$ws.Cells.Item(307, 7).Value = "Betarraga"
$ws.Cells.Item(307, 8).Value = "Sin especificar"
$ws.Cells.Item(307, 9).Value = "Primera"
$ws.Cells.Item(307, 10).Value = 3000
$ws.Cells.Item(307, 11).Value = 700
$ws.Cells.Item(307, 12).Value = 700
$ws.Cells.Item(307, 13).Value = 700
$ws.Cells.Item(307, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(307, 15).Value = "Región del Maule"
$ws.Cells.Item(307, 16).Value = 140
$ws.Cells.Item(307, 17).Value = 5
$ws.Cells.Item(307, 18).Value = "Hortaliza"
